$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (column type definitions) updates
$ws.Range("D1").Value = "varchar(50)"
$ws.Range("F1").Value = "varchar(6)"
$ws.Range("G1").Value = "int(20)"

# Row 2 (column labels) updates
$ws.Range("A2").Value = "ID (입력x)"
$ws.Range("B2").Value = "법인코드"
$ws.Range("C2").Value = "코스트센터코드"
$ws.Range("D2").Value = "version코드"
$ws.Range("E2").Value = "계정코드"
$ws.Range("F2").Value = "년월 ex) 200001"
$ws.Range("G2").Value = "금액"
